$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.8397313451970363
$ws.Range("C2").Value = 1.177508807361574
$ws.Range("D2").Value = 1.356451150188942
$ws.Range("E2").Value = 1.472823318998768
$ws.Range("B3").Value = 0.7676640602148203
$ws.Range("C3").Value = 1.087433508580475
$ws.Range("D3").Value = 1.276233810667159
$ws.Range("E3").Value = 1.409469658895177
$ws.Range("B4").Value = 0.4743663168598267
$ws.Range("C4").Value = 0.7013843045742871
$ws.Range("D4").Value = 0.9202580287282726
$ws.Range("E4").Value = 1.082814381615641
$ws.Range("B5").Value = 0.364676707402917
$ws.Range("C5").Value = 0.4306945007720887
$ws.Range("D5").Value = 0.4247129094724361
$ws.Range("E5").Value = 0.3980077702589124
$ws.Range("B6").Value = 0.2923816395688805
$ws.Range("C6").Value = 0.3093899654033974
$ws.Range("D6").Value = 0.2837220892666605
$ws.Range("E6").Value = 0.2390097767754088
$ws.Range("B7").Value = 0.328609636718072
$ws.Range("C7").Value = 0.3842135572168609
$ws.Range("D7").Value = 0.401576885974789
$ws.Range("E7").Value = 0.3706086512129885
$ws.Range("B8").Value = 0.1935278923803382
$ws.Range("C8").Value = 0.2819893251920724
$ws.Range("D8").Value = 0.3038021567545686
$ws.Range("E8").Value = 0.3085576887316853
$ws.Range("B9").Value = 0.5673291140965486
$ws.Range("C9").Value = 0.6558191130571009
$ws.Range("D9").Value = 0.5641100114228719
$ws.Range("E9").Value = 0.473483256811177
$ws.Range("B10").Value = 0.4268077567395996
$ws.Range("C10").Value = 0.7053690129767466
$ws.Range("D10").Value = 1.011569275504135
$ws.Range("E10").Value = 1.275327798525587
$ws.Range("B11").Value = 0.3745725955835346
$ws.Range("C11").Value = 0.6115985128189595
$ws.Range("D11").Value = 0.9020829115585566
$ws.Range("E11").Value = 1.162935609077621
$ws.Range("B12").Value = 0.09205845841638884
$ws.Range("C12").Value = 0.1751609156095072
$ws.Range("D12").Value = 0.391530093490966
$ws.Range("E12").Value = 0.614786416514029
$ws.Range("B13").Value = 0.3250309684184554
$ws.Range("C13").Value = 0.5463631393984159
$ws.Range("D13").Value = 0.8449974102226457
$ws.Range("E13").Value = 1.115428100200711